$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before F; the old "District" data in F shifts right to G.
$ws.Columns("F").Insert()

# New column F is the "Address" column. Header + a handful of populated rows.
$ws.Range("F2").Value = "Address"
$ws.Range("F6").Value = "Sri Gurushanteshwara High Sschool Godabanahal"
$ws.Range("F7").Value = "Govt High SchoolChavallihalliGollarahatti"
$ws.Range("F25").Value = "V B H S T R Nagar,Challakere"
$ws.Range("F27").Value = "GHS Valluru"
$ws.Range("F28").Value = "S Y B R High School Haikal"
$ws.Range("F31").Value = "Govt High School Badavanahalli"
$ws.Range("F36").Value = "Adarsha Vidyalaya(RMSA)"
$ws.Range("F41").Value = "Adarsha Vidyalaya Challakere"
$ws.Range("F43").Value = "Little Flower Girls High School (Aided)Hospet"
$ws.Range("F47").Value = "P M High SchoolAnkola"
$ws.Range("F50").Value = "Adarsha Vidyalaya Challakere"
